$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 1172.2142
$ws.Range("I6").Value = 1023.3333
$ws.Range("J6").Value = 1440.2
$ws.Range("K6").Value = 3069.9999
$ws.Range("L6").Value = 4320.6
$ws.Range("M6").Value = -2957.9999
$ws.Range("N6").Value = -4544.6
$ws.Range("H41").Value = 360.44
$ws.Range("I41").Value = 458.13333
$ws.Range("J41").Value = 213.9
$ws.Range("K41").Value = 458.13333
$ws.Range("L41").Value = 213.9
$ws.Range("M41").Value = -18.13333
$ws.Range("N41").Value = -1093.9
$ws.Range("H86").Value = 7257.143
$ws.Range("I86").Value = 5975
$ws.Range("K86").Value = 5975
$ws.Range("M86").Value = -4852
$ws.Range("H89").Value = 7257.143
$ws.Range("I89").Value = 5975
$ws.Range("K89").Value = 29875
$ws.Range("M89").Value = -24259
$ws.Range("H132").Value = 1135.5098
$ws.Range("I132").Value = 968.65
$ws.Range("J132").Value = 1742.2727
$ws.Range("K132").Value = 2905.95
$ws.Range("L132").Value = 5226.8181
$ws.Range("M132").Value = -375.9499999999998
$ws.Range("N132").Value = -10286.8181
$ws.Range("H133").Value = 118106
$ws.Range("J133").Value = 118106
$ws.Range("L133").Value = 118106
$ws.Range("N133").Value = -128226
$ws.Range("H138").Value = 3496.805
$ws.Range("J138").Value = 3844.6
$ws.Range("L138").Value = 11533.8
$ws.Range("N138").Value = -21813.8
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H19").Value = 1901.1428
$ws.Range("I19").Value = 2461.6
$ws.Range("K19").Value = 2461.6
$ws.Range("M19").Value = -2232.6
$ws.Range("H45").Value = 1169.5
$ws.Range("I45").Value = 1169.5
$ws.Range("K45").Value = 1169.5
$ws.Range("M45").Value = -792.5
$ws.Range("H61").Value = 8142.7407
$ws.Range("I61").Value = 6179.1333
$ws.Range("J61").Value = 10597.25
$ws.Range("K61").Value = 6179.1333
$ws.Range("L61").Value = 10597.25
$ws.Range("M61").Value = -5967.1333
$ws.Range("N61").Value = -11021.25
$ws.Range("H97").Value = 6174822.5
$ws.Range("I97").Value = 6174822.5
$ws.Range("K97").Value = 6174822.5
$ws.Range("M97").Value = -6174326.5
$ws.Range("H122").Value = 45637.332
$ws.Range("I122").Value = 4202.4
$ws.Range("J122").Value = 252812
$ws.Range("K122").Value = 12607.2
$ws.Range("L122").Value = 758436
$ws.Range("M122").Value = -10157.2
$ws.Range("N122").Value = -763336
$ws.Range("H125").Value = 49899
$ws.Range("J125").Value = 49899
$ws.Range("L125").Value = 49899
$ws.Range("N125").Value = -59739
$ws.Range("H132").Value = 9519.9
$ws.Range("I132").Value = 7281.0454
$ws.Range("K132").Value = 21843.1362
$ws.Range("M132").Value = -19313.1362
$ws.Range("H136").Value = 8142.7407
$ws.Range("I136").Value = 6179.1333
$ws.Range("J136").Value = 10597.25
$ws.Range("K136").Value = 18537.3999
$ws.Range("L136").Value = 31791.75
$ws.Range("M136").Value = -15987.3999
$ws.Range("N136").Value = -36891.75
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 44999.293
$ws.Range("I20").Value = 55367.105
$ws.Range("J20").Value = 5601.6
$ws.Range("K20").Value = 55367.105
$ws.Range("L20").Value = 5601.6
$ws.Range("M20").Value = -55120.105
$ws.Range("N20").Value = -6095.6
$ws.Range("H22").Value = 2136.9167
$ws.Range("I22").Value = 1374.1
$ws.Range("J22").Value = 5951
$ws.Range("K22").Value = 1374.1
$ws.Range("L22").Value = 5951
$ws.Range("M22").Value = -1201.1
$ws.Range("N22").Value = -6297
$ws.Range("H134").Value = 5087.491
$ws.Range("I134").Value = 3847.4443
$ws.Range("K134").Value = 11542.3329
$ws.Range("M134").Value = -9007.332900000001
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H17").Value = 33624.75
$ws.Range("I17").Value = 33624.75
$ws.Range("K17").Value = 33624.75
$ws.Range("M17").Value = -33450.75
$ws.Range("H31").Value = 27033836
$ws.Range("I31").Value = 111115530
$ws.Range("K31").Value = 111115530
$ws.Range("M31").Value = -111115235
$ws.Range("H34").Value = 27033836
$ws.Range("I34").Value = 111115530
$ws.Range("K34").Value = 111115530
$ws.Range("M34").Value = -111115328
$ws.Range("H41").Value = 38332.668
$ws.Range("J41").Value = 54999
$ws.Range("L41").Value = 54999
$ws.Range("N41").Value = -55855
$ws.Range("H50").Value = 59378.8
$ws.Range("J50").Value = 59378.8
$ws.Range("L50").Value = 59378.8
$ws.Range("N50").Value = -60628.8
$ws.Range("H51").Value = 40879.875
$ws.Range("J51").Value = 93363
$ws.Range("L51").Value = 93363
$ws.Range("N51").Value = -94835
$ws.Range("H58").Value = 5003.5713
$ws.Range("I58").Value = 2482.4
$ws.Range("J58").Value = 7912.615
$ws.Range("K58").Value = 2482.4
$ws.Range("L58").Value = 7912.615
$ws.Range("M58").Value = -2279.4
$ws.Range("N58").Value = -8318.615
$ws.Range("H59").Value = 40666.332
$ws.Range("J59").Value = 62249.75
$ws.Range("L59").Value = 62249.75
$ws.Range("N59").Value = -64539.75
$ws.Range("H60").Value = 7364.846
$ws.Range("J60").Value = 29453
$ws.Range("L60").Value = 29453
$ws.Range("N60").Value = -30475
$ws.Range("H61").Value = 40879.875
$ws.Range("J61").Value = 93363
$ws.Range("L61").Value = 93363
$ws.Range("N61").Value = -94059
$ws.Range("H132").Value = 16706.334
$ws.Range("I132").Value = 813.09375
$ws.Range("K132").Value = 2439.28125
$ws.Range("M132").Value = 90.71875
$ws.Range("H136").Value = 5003.5713
$ws.Range("I136").Value = 2482.4
$ws.Range("J136").Value = 7912.615
$ws.Range("K136").Value = 7447.200000000001
$ws.Range("L136").Value = 23737.845
$ws.Range("M136").Value = -4897.200000000001
$ws.Range("N136").Value = -28837.845
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 27210882
$ws.Range("I4").Value = 39297812
$ws.Range("J4").Value = 4648614.5
$ws.Range("K4").Value = 117893436
$ws.Range("L4").Value = 13945843.5
$ws.Range("M4").Value = -117893324
$ws.Range("N4").Value = -13946067.5
$ws.Range("H8").Value = 246.54546
$ws.Range("I8").Value = 246.54546
$ws.Range("K8").Value = 739.6363799999999
$ws.Range("M8").Value = -600.6363799999999
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("L16").ClearContents()
$ws.Range("M16").ClearContents()
$ws.Range("N16").Value = 0
$ws.Range("H33").Value = 242.88889
$ws.Range("I33").Value = 19.75
$ws.Range("J33").Value = 421.4
$ws.Range("K33").Value = 118.5
$ws.Range("L33").Value = 2528.4
$ws.Range("M33").Value = 164.5
$ws.Range("N33").Value = -3094.4
$ws.Range("H107").Value = 1006.4138
$ws.Range("J107").Value = 1695.75
$ws.Range("L107").Value = 5087.25
$ws.Range("N107").Value = -8927.25
$ws.Range("H122").Value = 58828610
$ws.Range("J122").Value = 992.1429000000001
$ws.Range("L122").Value = 8929.286100000001
$ws.Range("N122").Value = -13829.2861
$ws.Range("H134").Value = 71428840
$ws.Range("I134").Value = 71428840
$ws.Range("K134").Value = 214286520
$ws.Range("M134").Value = -214281450
$ws.Range("H136").Value = 2676.3333
$ws.Range("I136").Value = 2676.3333
$ws.Range("K136").Value = 8028.999899999999
$ws.Range("M136").Value = -2928.999899999999
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H33").Value = 24999
$ws.Range("J33").Value = 24999
$ws.Range("L33").Value = 24999
$ws.Range("N33").Value = -25503
$ws.Range("H80").Value = 7032.1665
$ws.Range("I80").Value = 6522.6
$ws.Range("J80").Value = 7396.143
$ws.Range("K80").Value = 6522.6
$ws.Range("L80").Value = 7396.143
$ws.Range("M80").Value = -5524.6
$ws.Range("N80").Value = -9392.143
$ws.Range("H83").Value = 7032.1665
$ws.Range("I83").Value = 6522.6
$ws.Range("J83").Value = 7396.143
$ws.Range("K83").Value = 32613
$ws.Range("L83").Value = 36980.715
$ws.Range("M83").Value = -27621
$ws.Range("N83").Value = -46964.715
$ws.Range("H132").Value = 4837.3716
$ws.Range("I132").Value = 3755
$ws.Range("K132").Value = 11265
$ws.Range("M132").Value = -8735
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 775.86664
$ws.Range("I9").Value = 99.09090999999999
$ws.Range("J9").Value = 2637
$ws.Range("K9").Value = 99.09090999999999
$ws.Range("L9").Value = 2637
$ws.Range("M9").Value = 124.90909
$ws.Range("N9").Value = -3085
$ws.Range("H26").Value = 6750
$ws.Range("J26").Value = 6750
$ws.Range("L26").Value = 6750
$ws.Range("N26").Value = -7340
$ws.Range("H46").Value = 6473.7144
$ws.Range("I46").Value = 1318.6666
$ws.Range("K46").Value = 1318.6666
$ws.Range("M46").Value = -1130.6666
$ws.Range("I68").Value = 2766.6667
$ws.Range("J68").Value = 7333.3335
$ws.Range("K68").Value = 2766.6667
$ws.Range("L68").Value = 7333.3335
$ws.Range("M68").Value = -2017.6667
$ws.Range("N68").Value = -8831.333500000001
$ws.Range("I71").Value = 2766.6667
$ws.Range("J71").Value = 7333.3335
$ws.Range("K71").Value = 13833.3335
$ws.Range("L71").Value = 36666.6675
$ws.Range("M71").Value = -10089.3335
$ws.Range("N71").Value = -44154.6675
$ws.Range("H132").Value = 5973.4
$ws.Range("I132").Value = 5357.1113
$ws.Range("J132").Value = 7046.9355
$ws.Range("K132").Value = 16071.3339
$ws.Range("L132").Value = 21140.8065
$ws.Range("M132").Value = -13541.3339
$ws.Range("N132").Value = -26200.8065
Write-Output "Updated market-price data across sheets."
